# Automatic update of files.
#
# Rows 12/13 and rows 14/15 each swap their observation data, i.e. the
# record that used to live in row 13 now lives in row 12 (and vice versa),
# and the record that used to live in row 15 now lives in row 14 (and vice
# versa). Row position therefore stays the same, but the content of each
# pair of rows trades places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Range, [string]$Text)
    if ([string]::IsNullOrEmpty($Text)) {
        $Range.Value2 = ""
    } elseif ($Text -match '^-?\d+(\.\d+)?$') {
        # Numeric-looking text must be forced to stay text (column I etc.)
        $Range.Value = "'" + $Text
    } else {
        $Range.Value2 = $Text
    }
}

# ---------------------------------------------------------------
# Row 12 <- values that used to belong to row 13 (Järpe / hazel grouse)
# ---------------------------------------------------------------
$ws.Range("A12").Value2 = 131196449
$ws.Range("B12").Value2 = 57064
$ws.Range("E12").Value2 = 102612
Set-TextCell $ws.Range("F12") "Järpe"
Set-TextCell $ws.Range("G12") "Tetrastes bonasia"
Set-TextCell $ws.Range("H12") "(Linnaeus, 1758)"
Set-TextCell $ws.Range("I12") "2"
Set-TextCell $ws.Range("M12") "födosökande"
Set-TextCell $ws.Range("N12") "observerad"
$ws.Range("Q12").Value2 = 500203
$ws.Range("R12").Value2 = 7016330
Set-TextCell $ws.Range("AC12") "Synobservation av 2 st födosökande järpar."

# ---------------------------------------------------------------
# Row 13 <- values that used to belong to row 12 (Garnlav / lichen)
# ---------------------------------------------------------------
$ws.Range("A13").Value2 = 131196451
$ws.Range("B13").Value2 = 79244
$ws.Range("E13").Value2 = 6425
Set-TextCell $ws.Range("F13") "Garnlav"
Set-TextCell $ws.Range("G13") "Alectoria sarmentosa"
Set-TextCell $ws.Range("H13") "(Ach.) Ach."
Set-TextCell $ws.Range("I13") ""
Set-TextCell $ws.Range("M13") ""
Set-TextCell $ws.Range("N13") ""
$ws.Range("Q13").Value2 = 500318
$ws.Range("R13").Value2 = 7016201
Set-TextCell $ws.Range("AC13") ""

# ---------------------------------------------------------------
# Row 14 <- values that used to belong to row 15 (Talltita / willow tit)
# ---------------------------------------------------------------
$ws.Range("A14").Value2 = 131196447
$ws.Range("B14").Value2 = 58043
$ws.Range("E14").Value2 = 103021
Set-TextCell $ws.Range("F14") "Talltita"
Set-TextCell $ws.Range("G14") "Poecile montanus"
Set-TextCell $ws.Range("H14") "(Conrad von Baldenstein, 1827)"
Set-TextCell $ws.Range("I14") "1"
Set-TextCell $ws.Range("M14") "förbiflygande"
Set-TextCell $ws.Range("N14") "observerad"
$ws.Range("Q14").Value2 = 500269
$ws.Range("R14").Value2 = 7016195
Set-TextCell $ws.Range("AC14") "Synobservation av 1 st talltita."

# ---------------------------------------------------------------
# Row 15 <- values that used to belong to row 14 (Garnlav / lichen)
# ---------------------------------------------------------------
$ws.Range("A15").Value2 = 131196452
$ws.Range("B15").Value2 = 79244
$ws.Range("E15").Value2 = 6425
Set-TextCell $ws.Range("F15") "Garnlav"
Set-TextCell $ws.Range("G15") "Alectoria sarmentosa"
Set-TextCell $ws.Range("H15") "(Ach.) Ach."
Set-TextCell $ws.Range("I15") ""
Set-TextCell $ws.Range("M15") ""
Set-TextCell $ws.Range("N15") ""
$ws.Range("Q15").Value2 = 500345
$ws.Range("R15").Value2 = 7016371
Set-TextCell $ws.Range("AC15") ""
